$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the site headers (A1:E1) -------------------------------
$ws.Range("A1").Value = "CAMPUS TECNOLÓGICO LOCAL SAN CARLOS"
$ws.Range("B1").Value = "CAMPUS TECNOLÓGICO LOCAL SAN JOSÉ"
$ws.Range("C1").Value = "CENTRO ACADÉMICO DE LIMÓN"
$ws.Range("D1").Value = "CAMPUS TECNOLÓGICO CENTRAL CARTAGO"
$ws.Range("E1").Value = "CENTRO ACADÉMICO DE ALAJUELA"

# --- 2. Strip the old "applyFill" cell style from every data cell ------
# (they go back to the default/Normal style). Only touch cells that
# already hold data so no new blank cells get materialised.
$ws.Range("A2:E3").Style = "Normal"
$ws.Range("A4:D4").Style = "Normal"
$ws.Range("A5:A8").Style = "Normal"
$ws.Range("D5:D21").Style = "Normal"

# --- 3. Add the new underlined, empty interface cell at C20 ------------
$ws.Range("C20").Font.Underline = $true
